$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Hope that this way of conceptuali" + bookmark + "zation could..."
#    -> merge into a single run and drop the _GoBack bookmark that used
#    to split it. (Paragraph 13: "Hope that this way of ...")
# ---------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$hopePara = $d.Paragraphs(13)
$hopeFull = $hopePara.Range
$hopeTextOnly = $d.Range($hopeFull.Start, $hopeFull.End - 1)
$hopeTextOnly.Delete()
$hopeTextOnly.InsertAfter("Hope that this way of conceptualization could be an alternative way to help you understand and remember the difference between recall and precision.")

# ---------------------------------------------------------------------
# 2) Add two new reference paragraphs (with hyperlinks) after the
#    existing "For more details, see this Wikipedia link" paragraph
#    (paragraph 23), and move the _GoBack bookmark to the very end of
#    the new content.
# ---------------------------------------------------------------------
$refPara = $d.Paragraphs(23)
$refPara.Range.InsertParagraphAfter()

$kagglePara = $d.Paragraphs(24)
$kagglePara.Range.Text = "See this Kaggle discussion about precision, recall, f-measure and AUC curve"

$kaggleLinkRange = $kagglePara.Range.Duplicate
$kaggleLinkRange.Find.Execute("precision, recall, f-measure and AUC curve") | Out-Null
$d.Hyperlinks.Add($kaggleLinkRange, "https://www.kaggle.com/general/7517", "", "", "precision, recall, f-measure and AUC curve") | Out-Null

$kagglePara.Range.InsertParagraphAfter()

$dsPara = $d.Paragraphs(25)
$dsPara.Range.Text = "See this guide from data school on precision, recall and on ROC curves and AUC"

$dsLinkRange1 = $dsPara.Range.Duplicate
$dsLinkRange1.Find.Execute("precision, recall") | Out-Null
$d.Hyperlinks.Add($dsLinkRange1, "http://www.dataschool.io/simple-guide-to-confusion-matrix-terminology/", "", "", "precision, recall") | Out-Null

$dsLinkRange2 = $dsPara.Range.Duplicate
$dsLinkRange2.Find.Execute("ROC curves and AUC") | Out-Null
$d.Hyperlinks.Add($dsLinkRange2, "http://www.dataschool.io/roc-curves-and-auc-explained/", "", "", "ROC curves and AUC") | Out-Null

$endRange = $dsPara.Range.Duplicate
$endRange.Find.Execute("AUC") | Out-Null
$endRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $endRange) | Out-Null
